$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "73.142.50"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "4.046.62"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'583.07"
$ws.Range("E5").Value = "  +10.29%  "
$ws.Range("D6").Value = "'152.03"
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("D7").Value = "4.041.59"
$ws.Range("E7").Value = "  +1.05%  "
$ws.Range("D8").Value = "'0.691"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "'0.758"
$ws.Range("E10").Value = "  +1.93%  "
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").Value = "'53.81"
$ws.Range("E12").Value = "  +13.30%  "
$ws.Range("D13").Value = "'0.0000324"
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("D14").Value = "'11.09"
$ws.Range("E14").Value = "  +4.48%  "
$ws.Range("D15").Value = "4.694.54"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "4.048.71"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").Value = "'14.27"
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("E18").Value = "  +4.73%  "
$ws.Range("D19").Value = "'20.77"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("E20").Value = "  -0.36%  "
$ws.Range("D21").Value = "73.130.19"
$ws.Range("E21").Value = "  +1.91%  "
$ws.Range("D22").Value = "'442.19"
$ws.Range("E22").Value = "  +3.58%  "
$ws.Range("D23").Value = "'4.61"
$ws.Range("E23").Value = "  +10.61%  "
$ws.Range("D24").Value = "'97.30"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("D26").Value = "'14.57"
$ws.Range("E26").Value = "  +1.54%  "
$ws.Range("D27").Value = "'4.32"
$ws.Range("E27").Value = "  +20.17%  "
$ws.Range("D28").Value = "'11.57"
$ws.Range("E28").Value = "  +2.62%  "
$ws.Range("E29").Value = "  +2.54%  "
$ws.Range("E30").Value = "  +2.04%  "
$ws.Range("D31").Value = "'36.89"
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("D32").Value = "'7.92"
$ws.Range("E32").Value = "  +11.59%  "
$ws.Range("E33").Value = "  +4.24%  "
$ws.Range("E34").Value = "  +2.44%  "
$ws.Range("D35").Value = "'689.17"
$ws.Range("E35").Value = "  +1.89%  "
$ws.Range("D36").Value = "'48.61"
$ws.Range("E36").Value = "  +9.10%  "
$ws.Range("D37").Value = "'67.34"
$ws.Range("E37").Value = "  +2.75%  "
$ws.Range("E38").Value = "  +2.10%  "
$ws.Range("D39").Value = "0.0₃0883"
$ws.Range("E39").Value = "  +7.13%  "
$ws.Range("D40").Value = "'0.150"
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("D41").Value = "'11.32"
$ws.Range("E41").Value = "  +18.40%  "
$ws.Range("D42").Value = "'3.36"
$ws.Range("E42").Value = "  -1.16%  "
$ws.Range("D44").Value = "'3.34"
$ws.Range("E44").Value = "  +4.60%  "
$ws.Range("E45").Value = "  +2.11%  "
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("E47").Value = "  +0.99%  "
$ws.Range("D48").Value = "'2.75"
$ws.Range("E48").Value = "  +4.58%  "
$ws.Range("D49").Value = "'3.37"
$ws.Range("E49").Value = "  -1.42%  "
$ws.Range("E50").Value = "  +7.10%  "
$ws.Range("D51").Value = "'3.05"
$ws.Range("E51").Value = "  +2.49%  "
